# Add employment dates after each role title in the "Employment Details"
# section. For each target paragraph we append two new runs: one holding a
# single space, and one holding the "(<start> - <end>)" date range — matching
# the structure already used elsewhere in the document (e.g. "::" + " " +
# italic title).

$d = $word.ActiveDocument

# Map of the (unique, exact) role-title paragraph text to the date range
# that should be appended after it.
$roleDates = @{
    "Head of Portfolio Architecture and Engineering" = "(Mar 2023 - Present)"
    "Cloud Practice Director"                        = "(Apr 2020 - Mar 2023)"
    "Technical Principal (CloudOps)"                 = "(Aug 2019 - Apr 2020)"
    "Azure Practice Lead"                             = "(Jul 2018 - Aug 2019)"
    "Windows Cloud Practice Lead - Bashton Ltd"      = "(Jul 2016 - Jul 2018)"
}

foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text
    # Paragraph text includes the trailing paragraph mark; strip it off so
    # we can compare against the expected title text exactly.
    $trimmed = $text.TrimEnd([char]13, [char]7)
    $title = $trimmed
    if ($title.StartsWith(":: ")) {
        $title = $title.Substring(3)
    }

    if ($roleDates.ContainsKey($title)) {
        $dateText = $roleDates[$title]

        # Range covering the paragraph text without its trailing mark, so
        # InsertAfter appends inside the paragraph (before the pilcrow).
        $pRange = $p.Range
        $insertPoint = $d.Range($pRange.Start, $pRange.End - 1)
        $insertPoint.InsertAfter(" ")

        $insertPoint2 = $d.Range($pRange.Start, $pRange.End - 1)
        $insertPoint2.InsertAfter($dateText)
    }
}
